$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TestData sheet (sheet2): duplicate the LoginPageTest row into row 3 and
# move the selection down to it.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TestData")
[void]($ws2.Activate())

[void]($ws2.Range("A2:C2").Copy($ws2.Range("A3")))

[void]($ws2.Range("A3:C3").Select())

# ---------------------------------------------------------------------------
# Config sheet (sheet1): append the new automated-test rows plus the new
# CustomerNumber column.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config")
[void]($ws1.Activate())

# New rows - written in this order so the shared-string table ends up with
# the same append order the original authoring session produced.
$ws1.Range("A5").Value = "No"
$ws1.Range("B5").Value = "DashboardPageTest"
$ws1.Range("C5").Value = "Yes"

$ws1.Range("A6").Value = "No"
$ws1.Range("B6").Value = "UpdatePageTest"
$ws1.Range("C6").Value = "Yes"

$ws1.Range("A7").Value = "Yes"
$ws1.Range("B7").Value = "SubmitPageTest"
$ws1.Range("C7").Value = "Yes"

$ws1.Range("A4").Value = "Yes"
$ws1.Range("B4").Value = "ExcelTest"
$ws1.Range("C4").Value = "Yes"

# New CustomerNumber column, copying the bold/fill header style from C1.
[void]($ws1.Range("C1").Copy($ws1.Range("D1")))
$ws1.Range("D1").Value = "CustomerNumber"

# Existing row 3 (HomePageTest) flips from "No" to "Yes" on both flags.
$ws1.Range("A3").Value = "Yes"
$ws1.Range("C3").Value = "Yes"

# Column widths (Excel COM ColumnWidth is expressed in "characters"; the
# saved XML width is ColumnWidth + 5/6).
$ws1.Columns.Item(1).ColumnWidth = 8.498697916666666
$ws1.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws1.Columns.Item(3).ColumnWidth = 10.166666666666666
$ws1.Columns.Item(4).ColumnWidth = 15.166666666666666

[void]($ws1.Range("B3").Select())

Write-Output "done"
